$d = $word.ActiveDocument

# 1. Collapse the "Efficient access to elements" bullet into a single run
#    (removes the proofErr gramStart/gramEnd markers around "provide" and
#    merges the three runs back into one run with the same full text).
$find = $d.Content.Find
$effText = "Efficient access to elements: Arrays provide direct and efficient access to any element in the collection. Accessing an element in an array is an O(1) operation, meaning that the time required to access an element is constant and does not depend on the size of the array."
$find.Execute($effText, $true, $false, $false, $false, $false, $true, 1, $false, $effText, 2) | Out-Null

# 2. Remove the old "Memory efficiency" bullet (it will be recreated from
#    the "Fast data retrieval" bullet below), the "Versatility" bullet, and
#    the "Compatibility with hardware" bullet entirely (paragraph and all).
#    Walk bottom-up so earlier paragraph indices stay valid while deleting.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Memory efficiency: Arrays are a memory-efficient way of storing data.")) {
        $p.Range.Delete()
    }
    elseif ($t.StartsWith("Versatility: Arrays can be used to store a wide range of data types")) {
        $p.Range.Delete()
    }
    elseif ($t.StartsWith("Compatibility with hardware: The array data structure is compatible")) {
        $p.Range.Delete()
    }
}

# 3. Replace the "Fast data retrieval" bullet text with the new
#    "Memory efficiency" text.
$find = $d.Content.Find
$oldFast = "Fast data retrieval: Arrays allow for fast data retrieval because the data is stored in contiguous memory locations. This means that the data can be accessed quickly and efficiently without the need for complex data structures or algorithms."
$newMem = "Memory efficiency: Arrays are a memory-efficient way of storing data. Because the elements of an array are stored in contiguous memory locations, the size of the array is known at compile time. This means that memory can be allocated for the entire array in one block, reducing memory fragmentation."
$find.Execute($oldFast, $true, $false, $false, $false, $false, $true, 1, $false, $newMem, 2) | Out-Null
